$wb = $excel.ActiveWorkbook

# --- Sheet references -------------------------------------------------
$wsBusIndex = $wb.Worksheets.Item("Bus index")
$wsBusConn  = $wb.Worksheets.Item("Bus connections")
$wsNodeIdx  = $wb.Worksheets.Item("Node index")
$wsGenData  = $wb.Worksheets.Item("Generator data")

# --- "Bus index" sheet: D2 content + formatting ------------------------
# D2 changes text (new shared string) and picks up the "orange" look
# (fontId 4 / fillId 3, like the Node index!E2:E5 cells) but left aligned
# instead of centered.
$wsNodeIdx.Range("E2").Copy() | Out-Null
$wsBusIndex.Range("D2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsBusIndex.Range("D2").Value = "Bubble HUN Wind 2014-2045_0910refyr.csv"
$wsBusIndex.Range("D2").HorizontalAlignment = -4131      # xlLeft

# --- "Bus connections" sheet: reactance formulas + MW limits -----------
$wsBusConn.Range("C2").Formula = "=(PI()/12)/(E2/100)*D2"
$wsBusConn.Range("C3").Formula = "=(PI()/12)/(E3/100)*D3"

$wsBusConn.Range("E4").Value = 2700
$wsBusConn.Range("C4").Formula = "=(PI()/6)/(E4/100)*D4"

$wsBusConn.Range("E5").Value = 1800
$wsBusConn.Range("C5").Formula = "=(PI()/3)/(E5/100)*D5"

$wsBusConn.Range("E6").Value = 4800
$wsBusConn.Range("E8").Value = 4600

# --- Selections / active sheet -----------------------------------------
# Touch the non-active sheets first so their selection gets recorded
# without leaving them as the active tab, then finish on "Bus index" so
# it ends up the visible/active sheet.
$wsBusConn.Range("E4").Select() | Out-Null
$wsNodeIdx.Range("B38").Select() | Out-Null
$wsBusIndex.Range("D2").Select() | Out-Null
